$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing periods D:K to E:L
$ws.Columns("D:D").Insert()

# Copy number formats from column E (the old column D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest reporting period data
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 388379000
$ws.Range("D9").Value2 = 321369000
$ws.Range("D10").Value2 = 67010000
$ws.Range("D12").Value2 = 2326000
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 431000
$ws.Range("D15").Value2 = 21704000
$ws.Range("D17").Value2 = 357988000
$ws.Range("D18").Value2 = 30391000
$ws.Range("D20").Value2 = 7904000
$ws.Range("D21").Value2 = 60430000
$ws.Range("D22").Value2 = 2674000
$ws.Range("D23").Value2 = 35621000
$ws.Range("D24").Value2 = 11715000
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 23906000
$ws.Range("D27").Value2 = 23352000
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -7904000
$ws.Range("D33").Value2 = 23352000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 23352000
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 3777000
$ws.Range("D42").Value2 = 22707000
$ws.Range("D43").Value2 = 37076000
$ws.Range("D44").Value2 = 21117000
$ws.Range("D45").Value2 = 12805000
$ws.Range("D46").Value2 = 97482000
$ws.Range("D47").Value2 = 34409000
$ws.Range("D48").Value2 = 223175000
$ws.Range("D49").Value2 = 23586000
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 20542000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 399194000
$ws.Range("D57").Value2 = 33202000
$ws.Range("D58").Value2 = 10134000
$ws.Range("D59").Value2 = 34477000
$ws.Range("D60").Value2 = 77813000
$ws.Range("D61").Value2 = 66690000
$ws.Range("D62").Value2 = 52157000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 200548000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 221097000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 198646000
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 23352000
$ws.Range("D83").Value2 = 22135000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 53085000
$ws.Range("D91").Value2 = -23011000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -13659000
$ws.Range("D96").Value2 = -15675000
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -32548000
$ws.Range("D101").Value2 = -449000
$ws.Range("D102").Value2 = 6429000

# Restate historical Accounts Payable (row 57) and Other Current Liabilities (row 59)
$ws.Range("E57").Value2 = 35217000
$ws.Range("F57").Value2 = 30042000
$ws.Range("G57").Value2 = 26298000
$ws.Range("H57").Value2 = 34833000
$ws.Range("I57").Value2 = 45821000
$ws.Range("J57").Value2 = 89764000
$ws.Range("E59").Value2 = 32755000
$ws.Range("F59").Value2 = 34299000
$ws.Range("G59").Value2 = 39120000
$ws.Range("H59").Value2 = 44171000
$ws.Range("I59").Value2 = 39093000
$ws.Range("J59").Value2 = 63039000
